$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the data table with a new "2021" column (column R), mirroring the
# existing formatting of column Q (the previous last column, "2020").
# Copying Q4:Q14 -> R4:R14 brings over number formats/fonts/borders so the
# new column visually matches its neighbours, then we overwrite just the
# cells whose 2021 figure actually differs from the 2020 one that got
# copied in (several indicators repeat their 2020 value for 2021).
$ws.Range("Q4:Q14").Copy($ws.Range("R4:R14"))

$ws.Range("R4").Value = 2021
$ws.Range("R5").Value = 99.4
$ws.Range("R6").Value = 98.1
$ws.Range("R9").Value = 99.1
$ws.Range("R10").Value = 99.3
$ws.Range("R12").Value = 99.3

# Match the workbook's saved selection/active cell.
$ws.Range("U4").Select()
